# Colocando header nos gráficos
$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheets 1-4 share the same "Fonte/Tecnologia" table layout (rows 1-12):
#   - Add a header label in A1 (using the same style as the other header
#     cells in row 1)
#   - Strip the header style from A2:A12 (they become plain/unstyled cells)
#   - Fix a few accented labels
# -------------------------------------------------------------------------
$sourceSheets = @(1, 2, 3, 4)

foreach ($idx in $sourceSheets) {
    $ws = $wb.Worksheets.Item($idx)

    # Add header cell A1
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

    # Correct accented labels
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."

    # Remove the (bold/bordered) header style from the row labels A2:A12
    $ws.Range("A2:A12").ClearFormats()
}

$excel.CutCopyMode = $false

# -------------------------------------------------------------------------
# Sheet 5 "Emissoes Totais (MtCO2eq)"
#   - Add header label in A1 ("Período")
#   - Strip header style from A2:A3
#   - Fix accented labels
#   - Remove row 4 ("Teto") entirely
# -------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A1").Value = "Período"
$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Range("A2:A3").ClearFormats()

$ws5.Rows.Item(4).Delete()

# -------------------------------------------------------------------------
# Sheet 6 "Custo Total (bilhões de R$)"
#   - Add header label in A1 ("Tipo Expansão")
#   - B1 header text changes from "Custo" to "2015"
#   - Strip header style from A2:A3
#   - Fix accented labels
#   - Update values in B2/B3
# -------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws1 = $wb.Worksheets.Item(1)

# Header text "Custo" -> "2015" (use a leading apostrophe so the numeric
# looking text stays a text value, then re-apply the clean header format
# from sheet 1's B1 so the style index matches, without a quote-prefix flag)
$ws6.Range("B1").Value = "'2015"
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 548

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99

$ws6.Range("A2:A3").ClearFormats()
